# Apply weekly update to the "Zapallo italiano" hortaliza sheet.
# The edit permutes the per-row data (Fecha, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg o Unidades) among rows
# 3, 4, 5, 6, 7, 10 and 11, while leaving the remaining (identifying) columns intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows.
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Capture the current ("before") values for every row that participates in the
# permutation so we can reassign them without clobbering data we still need.
$rows = @(3, 4, 5, 6, 7, 10, 11)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# new row -> source row (where its new content currently lives)
$mapping = @{
    3  = 6
    4  = 10
    5  = 7
    6  = 11
    7  = 5
    10 = 4
    11 = 3
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
